# Update the "Metadata" sheet: IG publication Date moved forward a day.
$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-07-04T07:50:29+01:00"

# Update the "Elements" sheet: the Extension.value[x] binding strength
# moved from "required" to "extensible", and the bound ValueSet URL was
# repointed to the new (https, immunizationIG-path) location.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("X6").Value = "extensible"
$elements.Range("Z6").Value = "https://nphcda.gov.ng/immunizationIG/ValueSet/nigeria-vaccine-contraindication"

# Widen column Z ("Binding Value Set") so the longer URL still best-fits.
$elements.Columns.Item(26).ColumnWidth = 64
